$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = 60
$ws1.Range("F6").Value = 5326
$ws1.Range("F8").Value = 878
$ws1.Range("F9").Value = 127
$ws1.Range("F10").Value = 2367
$ws1.Range("G10").Value = 65
$ws1.Range("F12").Value = 49
$ws1.Range("F13").Value = 2218

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = 60
$ws4.Range("F6").Value = 5326
$ws4.Range("F10").Value = 878
$ws4.Range("F11").Value = 127
$ws4.Range("F12").Value = 2367
$ws4.Range("G12").Value = 65
$ws4.Range("F15").Value = 49
$ws4.Range("F16").Value = 2218
